# Auto-generated: refresh market-data columns (H-N) across all Leve profit sheets
# per scheduled-runner commit. Values below come from the authoritative diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 180.11111
$ws.Range("I31").Value = 180.11111
$ws.Range("K31").Value = 540.3333299999999
$ws.Range("M31").Value = -310.3333299999999
$ws.Range("H33").Value = 435.93332
$ws.Range("I33").Value = 445.64285
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 445.64285
$ws.Range("L33").Value = 300
$ws.Range("M33").Value = -216.64285
$ws.Range("N33").Value = -758
$ws.Range("H62").Value = 2687.5
$ws.Range("I62").Value = 1916.6666
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 1916.6666
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -1292.6666
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 2687.5
$ws.Range("I65").Value = 1916.6666
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 9583.333000000001
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -6463.333000000001
$ws.Range("N65").Value = -31240
$ws.Range("H92").Value = 335.66666
$ws.Range("I92").Value = 327.625
$ws.Range("K92").Value = 327.625
$ws.Range("M92").Value = 920.375
$ws.Range("H98").Value = 4520.2
$ws.Range("I98").Value = 4664.5
$ws.Range("J98").Value = 3221.5
$ws.Range("K98").Value = 4664.5
$ws.Range("L98").Value = 3221.5
$ws.Range("M98").Value = -3166.5
$ws.Range("N98").Value = -6217.5
$ws.Range("H122").Value = 4520.2
$ws.Range("I122").Value = 4664.5
$ws.Range("J122").Value = 3221.5
$ws.Range("K122").Value = 13993.5
$ws.Range("L122").Value = 9664.5
$ws.Range("M122").Value = -11543.5
$ws.Range("N122").Value = -14564.5
$ws.Range("H137").Value = 3825.7878
$ws.Range("I137").Value = 3087.52
$ws.Range("K137").Value = 9262.559999999999
$ws.Range("M137").Value = -6712.559999999999
$ws.Range("H138").Value = 3365.1516
$ws.Range("I138").Value = 3691
$ws.Range("K138").Value = 11073
$ws.Range("M138").Value = -5933
$ws.Range("H141").Value = 3681.1667
$ws.Range("I141").Value = 3617.6
$ws.Range("K141").Value = 10852.8
$ws.Range("M141").Value = -5672.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 867.2857
$ws.Range("I4").Value = 838.6
$ws.Range("K4").Value = 838.6
$ws.Range("M4").Value = -722.6
$ws.Range("H32").Value = 1905.75
$ws.Range("I32").Value = 1910
$ws.Range("K32").Value = 1910
$ws.Range("M32").Value = -1623
$ws.Range("H63").Value = 2364.75
$ws.Range("I63").Value = 2364.75
$ws.Range("K63").Value = 2364.75
$ws.Range("M63").Value = -1678.75
$ws.Range("H66").Value = 2364.75
$ws.Range("I66").Value = 2364.75
$ws.Range("K66").Value = 11823.75
$ws.Range("M66").Value = -8391.75
$ws.Range("H97").Value = 745.53845
$ws.Range("J97").Value = 165
$ws.Range("L97").Value = 165
$ws.Range("N97").Value = -1157

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4313.875
$ws.Range("I86").Value = 3503.6667
$ws.Range("K86").Value = 3503.6667
$ws.Range("M86").Value = -2380.6667
$ws.Range("H89").Value = 4313.875
$ws.Range("I89").Value = 3503.6667
$ws.Range("K89").Value = 17518.3335
$ws.Range("M89").Value = -11902.3335
$ws.Range("H99").Value = 2778.0667
$ws.Range("I99").Value = 1818.3
$ws.Range("J99").Value = 4697.6
$ws.Range("K99").Value = 1818.3
$ws.Range("L99").Value = 4697.6
$ws.Range("M99").Value = -320.3
$ws.Range("N99").Value = -7693.6
$ws.Range("H61").Value = 50000
$ws.Range("J61").Value = 50000
$ws.Range("L61").Value = 50000
$ws.Range("N61").Value = -50626
$ws.Range("H132").Value = 85119.2
$ws.Range("I132").Value = 32000
$ws.Range("J132").Value = 98399
$ws.Range("K132").Value = 32000
$ws.Range("L132").Value = 98399
$ws.Range("M132").Value = -26940
$ws.Range("N132").Value = -108519

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 261783.33
$ws.Range("I4").Value = 125175
$ws.Range("K4").Value = 125175
$ws.Range("M4").Value = -125063
$ws.Range("H10").Value = 2500
$ws.Range("I10").Value = 2500
$ws.Range("K10").Value = 2500
$ws.Range("M10").Value = -2361
$ws.Range("H31").Value = 2407651.2
$ws.Range("I31").Value = 4036.353
$ws.Range("K31").Value = 4036.353
$ws.Range("M31").Value = -3741.353
$ws.Range("H34").Value = 2407651.2
$ws.Range("I34").Value = 4036.353
$ws.Range("K34").Value = 4036.353
$ws.Range("M34").Value = -3834.353
$ws.Range("H86").Value = 10305.538
$ws.Range("I86").Value = 8897.5
$ws.Range("K86").Value = 8897.5
$ws.Range("M86").Value = -7774.5
$ws.Range("H89").Value = 10305.538
$ws.Range("I89").Value = 8897.5
$ws.Range("K89").Value = 44487.5
$ws.Range("M89").Value = -38871.5
$ws.Range("H94").Value = 941.4286
$ws.Range("I94").Value = 1188
$ws.Range("J94").Value = 900.3333
$ws.Range("K94").Value = 1188
$ws.Range("L94").Value = 900.3333
$ws.Range("M94").Value = -737
$ws.Range("N94").Value = -1802.3333
$ws.Range("H107").Value = 2500804.2
$ws.Range("I107").Value = 3846769.2
$ws.Range("K107").Value = 3846769.2
$ws.Range("M107").Value = -3844849.2
$ws.Range("H122").Value = 1285.0435
$ws.Range("I122").Value = 1416.7646
$ws.Range("K122").Value = 4250.293799999999
$ws.Range("M122").Value = -1800.293799999999
$ws.Range("H132").Value = 6293597.5
$ws.Range("I132").Value = 4204.5713
$ws.Range("J132").Value = 30307644
$ws.Range("K132").Value = 12613.7139
$ws.Range("L132").Value = 90922932
$ws.Range("M132").Value = -10083.7139
$ws.Range("N132").Value = -90927992
$ws.Range("H134").Value = 3663.1785
$ws.Range("I134").Value = 3259.7
$ws.Range("J134").Value = 4671.875
$ws.Range("K134").Value = 9779.099999999999
$ws.Range("L134").Value = 14015.625
$ws.Range("M134").Value = -7244.099999999999
$ws.Range("N134").Value = -19085.625
$ws.Range("H2").Value = 2000
$ws.Range("I2").Value = 2000
$ws.Range("K2").Value = 2000
$ws.Range("M2").Value = -1887
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1721.3334
$ws.Range("I92").Value = 1047
$ws.Range("J92").Value = 1914
$ws.Range("K92").Value = 3141
$ws.Range("L92").Value = 5742
$ws.Range("M92").Value = -1893
$ws.Range("N92").Value = -8238
$ws.Range("H139").Value = 1987.75
$ws.Range("I139").Value = 1097.5
$ws.Range("K139").Value = 3292.5
$ws.Range("M139").Value = 1847.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2006.8889
$ws.Range("I102").Value = 1558.1364
$ws.Range("J102").Value = 3981.4
$ws.Range("K102").Value = 1558.1364
$ws.Range("L102").Value = 3981.4
$ws.Range("M102").Value = 63.86359999999991
$ws.Range("N102").Value = -7225.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1168
$ws.Range("I55").Value = 1047.6666
$ws.Range("J55").Value = 1529
$ws.Range("K55").Value = 1047.6666
$ws.Range("L55").Value = 1529
$ws.Range("M55").Value = -874.6666
$ws.Range("N55").Value = -1875
$ws.Range("H93").Value = 2460.375
$ws.Range("I93").Value = 2526.1428
$ws.Range("K93").Value = 2526.1428
$ws.Range("M93").Value = -1278.1428
$ws.Range("H136").Value = 2171.7058
$ws.Range("I136").Value = 1327.6666
$ws.Range("J136").Value = 3121.25
$ws.Range("K136").Value = 3982.9998
$ws.Range("L136").Value = 9363.75
$ws.Range("M136").Value = -1432.9998
$ws.Range("N136").Value = -14463.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7023.75
$ws.Range("I62").Value = 5573.5
$ws.Range("J62").Value = 8474
$ws.Range("K62").Value = 5573.5
$ws.Range("L62").Value = 8474
$ws.Range("M62").Value = -4949.5
$ws.Range("N62").Value = -9722
$ws.Range("H65").Value = 7023.75
$ws.Range("I65").Value = 5573.5
$ws.Range("J65").Value = 8474
$ws.Range("K65").Value = 27867.5
$ws.Range("L65").Value = 42370
$ws.Range("M65").Value = -24747.5
$ws.Range("N65").Value = -48610
$ws.Range("H96").Value = 4145.75
$ws.Range("I96").Value = 2329.5
$ws.Range("K96").Value = 2329.5
$ws.Range("M96").Value = -956.5
$ws.Range("H122").Value = 11364719
$ws.Range("I122").Value = 1103.6842
$ws.Range("K122").Value = 3311.0526
$ws.Range("M122").Value = -861.0526
$ws.Range("H136").Value = 209281.3
$ws.Range("I136").Value = 3697.3125
$ws.Range("K136").Value = 11091.9375
$ws.Range("M136").Value = -8541.9375
